$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.618479490280151
$ws.Range("B1").Value = 2.526807546615601
$ws.Range("C1").Value = 2.759284496307373
$ws.Range("D1").Value = 3.091423511505127
$ws.Range("E1").Value = 3.405114889144897
